{"js": "// Revert \"second para added\": remove the second paragraph (\"Hey potatoes...\")\n// and restore the first paragraph's original run split, where the word\n// \"songwriting\" sits in its own run bracketed by spell-check proofErr marks.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load('items/text');\nawait context.sync();\n\n// Remove the extra \"Hey potatoes, im soooooooooo bored. wbu\" paragraph\n// that was appended by the commit this change reverts.\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf('Hey potatoes') !== -1) {\n    p.delete();\n  }\n}\nawait context.sync();\n\n// Re-split the bio paragraph's single run into three runs, with\n// <w:proofErr w:type=\"spellStart\"/>...<w:proofErr w:type=\"spellEnd\"/>\n// wrapped around the \"songwriting\" run, matching the pre-edit markup.\nconst bodyParagraphs = body.paragraphs;\nbodyParagraphs.load('items/text');\nawait context.sync();\n\nlet bioParagraph = null;\nfor (const p of bodyParagraphs.items) {\n  if (p.text.indexOf('songwriting') !== -1) {\n    bioParagraph = p;\n    break;\n  }\n}\nif (!bioParagraph) {\n  bioParagraph = bodyParagraphs.items[0];\n}\nconst bioRange = bioParagraph.getRange();\n\nconst ooxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n          '<w:body>' +\n            '<w:p>' +\n              '<w:r><w:t xml:space=\"preserve\">Taylor Alison Swift (born December 13, 1989) is an American singer-songwriter. Her narrative </w:t></w:r>' +\n              '<w:proofErr w:type=\"spellStart\"/>' +\n              '<w:r><w:t>songwriting</w:t></w:r>' +\n              '<w:proofErr w:type=\"spellEnd\"/>' +\n              '<w:r><w:t>, which often takes inspiration from her personal life, has received widespread critical praise and media coverage.</w:t></w:r>' +\n            '</w:p>' +\n          '</w:body>' +\n        '</w:document>' +\n      '</pkg:xmlData>' +\n    '</pkg:part>' +\n  '</pkg:package>';\n\nbioRange.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Revert \"second para added\": remove the second paragraph (\"Hey potatoes...\")\n# and restore the first paragraph's original run split, where the word\n# \"songwriting\" sits in its own run bracketed by spell-check proofErr marks.\n\n$d = $word.ActiveDocument\n\n# Remove the extra \"Hey potatoes, im soooooooooo bored. wbu\" paragraph that\n# was appended by the commit this change reverts (walk backwards since\n# deleting shifts later indices).\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $para = $d.Paragraphs.Item($i)\n    if ($para.Range.Text -like '*Hey potatoes*') {\n        $para.Range.Delete()\n    }\n}\n\n# Re-split the bio paragraph's single run into three runs, with\n# <w:proofErr w:type=\"spellStart\"/>...<w:proofErr w:type=\"spellEnd\"/>\n# wrapped around the \"songwriting\" run, matching the pre-edit markup.\n$bio = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $para = $d.Paragraphs.Item($i)\n    if ($para.Range.Text -like '*songwriting*') {\n        $bio = $para.Range\n        break\n    }\n}\nif ($null -eq $bio) {\n    $bio = $d.Paragraphs.Item(1).Range\n}\n$bio.MoveEnd(1, -1) | Out-Null   # exclude the trailing paragraph mark\n\n$xml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n            '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n                '<w:body>' +\n                    '<w:p>' +\n                        '<w:r><w:t xml:space=\"preserve\">Taylor Alison Swift (born December 13, 1989) is an American singer-songwriter. Her narrative </w:t></w:r>' +\n                        '<w:proofErr w:type=\"spellStart\"/>' +\n                        '<w:r><w:t>songwriting</w:t></w:r>' +\n                        '<w:proofErr w:type=\"spellEnd\"/>' +\n                        '<w:r><w:t>, which often takes inspiration from her personal life, has received widespread critical praise and media coverage.</w:t></w:r>' +\n                    '</w:p>' +\n                '</w:body>' +\n            '</w:document>' +\n        '</pkg:xmlData>' +\n    '</pkg:part>' +\n'</pkg:package>'\n\n$bio.InsertXML($xml) | Out-Null\n"}
